$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = 1583971200

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "2020-03-12"
$ws.Range("B17").Style = $ws.Range("B16").Style

$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "5295"
$ws.Range("C17").Style = $ws.Range("C16").Style

$ws.Range("D17").Value = "INNATURE"
$ws.Range("E17").Value = 0.43
$ws.Range("F17").Value = 0.43
$ws.Range("G17").Value = 0.38
$ws.Range("H17").Value = 0.39
$ws.Range("I17").Value = 1710700
